# Automatische test-sync: 2025-06-26 19:26:50
# Adds a new log row (row 9) to the "Logs" sheet, bumps the "Bestelling / Levering"
# tally on the "Dashboard" sheet, and extends the conditional-formatting ranges
# that covered the old used range (rows 2-8) so they also cover the new row 9.

$wb = $excel.ActiveWorkbook

# ---- Logs sheet: append the new mail-log entry as row 9 ----
$logs = $wb.Worksheets.Item("Logs")

$logs.Range("A9").Value = "Bestel je 3 rollen ducttape?"
$logs.Range("B9").Value = "MailMind Test <mailmind.test@zohomail.eu>"
$logs.Range("C9").Value = "He Johan,`nZou je 3 rollen ducttape kunnen bestellen?`nMarc`nSent using {0}"
$logs.Range("D9").Value = "Bestelling / Levering"
$logs.Range("E9").Value = "Beste Marc,`nBedankt voor je verzoek. Ik zal direct 3 rollen ducttape bestellen. Heb je een specifieke voorkeur voor het merk of kleur? Laat het me weten, zodat ik de bestelling nauwkeurig kan plaatsen.`nMet vriendelijke groet,`nJohan"
$logs.Range("F9").Value = "2025-06-26 19:26:42"
$logs.Range("G9").Value = "Ja"
$logs.Range("H9").Value = "Nee"
$logs.Range("I9").Value = "Ja"

# Extend the conditional formatting ranges from row 2-8 to row 2-9, matching the
# newly grown used range.
$logs.Range("D2:D8").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("D2:D9"))
$logs.Range("G2:G8").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("G2:G9"))
$logs.Range("H2:H8").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("H2:H9"))
$logs.Range("I2:I8").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("I2:I9"))

# ---- Dashboard sheet: bump the "Bestelling / Levering" count from 3 to 4 ----
$dashboard = $wb.Worksheets.Item("Dashboard")
$dashboard.Range("B2").Value = 4
